$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("routes.txt")

# Insert a new column before the existing "route_color" column (F),
# shifting route_color/route_text_color one column to the right.
$ws.Columns.Item(6).Insert()

# New header for the inserted "route_url" column (the insert already
# carried over the bold header style from the surrounding cells).
$ws.Range("F1").Value = "route_url"
$ws.Columns.Item(6).ColumnWidth = 15.67

# Make routes.txt the active sheet/tab and select cell F2, matching the
# saved view state of the edited workbook.
$ws.Activate()
[void]$ws.Range("F2").Select()
